# Atualização automática de preços de eletricidade
# Updates the single data row (row 2) of the Spot_PT sheet with the
# newly published day and hourly prices.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Day (date serial number)
$ws.Range("A2").Value = 45971

# Hourly prices (0h-1h .. 23h-24h)
$ws.Range("B2").Value  = 89.33
$ws.Range("C2").Value  = 78.7
$ws.Range("D2").Value  = 66.34999999999999
$ws.Range("E2").Value  = 59.93
$ws.Range("F2").Value  = 59.68
$ws.Range("G2").Value  = 66.87
$ws.Range("H2").Value  = 91.20999999999999
$ws.Range("I2").Value  = 110.06
$ws.Range("J2").Value  = 111.6
$ws.Range("K2").Value  = 72.66
$ws.Range("L2").Value  = 44.11
$ws.Range("M2").Value  = 35.68
$ws.Range("N2").Value  = 37.72
$ws.Range("O2").Value  = 31
$ws.Range("P2").Value  = 32.78
$ws.Range("Q2").Value  = 47.77
$ws.Range("R2").Value  = 70.15000000000001
$ws.Range("S2").Value  = 110.49
$ws.Range("T2").Value  = 131.57
$ws.Range("U2").Value  = 143.68
$ws.Range("V2").Value  = 136.88
$ws.Range("W2").Value  = 131.29
$ws.Range("X2").Value  = 107.02
$ws.Range("Y2").Value  = 93.97
$ws.Range("Z2").Value  = 81.69

# Daily average (Price_Daily_Avg) stays the same slot label
$ws.Range("AA2").Value = "20h-24h"

# Slot_4h_price
$ws.Range("AB2").Value = 117.29

# Slot_2h_frist / Slot_2h_frist_price
$ws.Range("AC2").Value = "18h-20h"
$ws.Range("AD2").Value = 137.62

# Slot_2h_second / Slot_2h_second_price
$ws.Range("AE2").Value = "20h-22h"
$ws.Range("AF2").Value = 134.08

# Slot_min_price
$ws.Range("AG2").Value = "1h-16h"
